$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 33: copy the formatting of the previous last data row (32) first,
# so the new row reuses the existing font/alignment/border styles instead of
# minting brand-new ones, then drop the bottom border (row 33 becomes the new
# last row so the border no longer needs to close the table underneath it).
$srcRow = $ws.Range("A32:B32")
$newRow = $ws.Range("A33:B33")
$srcRow.Copy()
$newRow.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$newRow.Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> none (xlLineStyleNone)

# Data for the new row
$ws.Range("A33").Value = 18589
$ws.Range("B33").Value = "沃游戏"

# Row height for the new row
$ws.Rows.Item(33).RowHeight = 16.5

# New selection location
$ws.Range("A2").Select() | Out-Null
